$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 704, pushing existing rows 704-820 down to 705-821.
$ws.Rows.Item(704).Insert()

# Populate the newly inserted row 704 with the new data record.
$ws.Range("A704").Value = 3
$ws.Range("B704").Value = "Femacal de La Calera"
$ws.Range("C704").Value = "Coquimbo"
$ws.Range("D704").Value = 45180
$ws.Range("E704").Value = 5
$ws.Range("F704").Value = 100112037
$ws.Range("G704").Value = "Cebollín"
$ws.Range("H704").Value = "Sin especificar"
$ws.Range("I704").Value = "Primera"
$ws.Range("J704").Value = 230
$ws.Range("K704").Value = 4000
$ws.Range("L704").Value = 4500
$ws.Range("M704").Value = 4239
$ws.Range("N704").Value = "$/paquete 36 unidades"
$ws.Range("O704").Value = "Provincia de Quillota"
$ws.Range("P704").Value = 118
$ws.Range("Q704").Value = 36
$ws.Range("R704").Value = "Hortaliza"
